# Add a new "Blood Amount Preset" row into the CategoryPresetSelection block
# (new row 8) and shift every subsequent row down by one. Then fix up the
# DOT damage-type multiplier defaults (now rows 9-13) to match the new
# screenshot-sourced values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 8 - pushes old rows 8..56 down to 9..57
$ws.Rows.Item(8).Insert()

# New row 8: CategoryPresetSelection / OptionBloodAmountPreset
$ws.Cells.Item(8, 1).Value = "CategoryPresetSelection"
$ws.Cells.Item(8, 2).Value = 50
$ws.Cells.Item(8, 3).Value = "OptionBloodAmountPreset"
$ws.Cells.Item(8, 4).Value = "string"
$ws.Cells.Item(8, 5).Value = """Default"""
$ws.Cells.Item(8, 6).Value = "Blood VFX intensity preset. Controls how much blood spurts from wounds. Very Low = minimal blood, Default = moderate blood, Extreme = lots of blood."

# Fix up the damage type multiplier defaults (old rows 8-12, now 9-13)
$ws.Cells.Item(9, 5).Value = "1.2f"    # OptionPierceMultiplier
$ws.Cells.Item(10, 5).Value = "0.8f"   # OptionSlashMultiplier
$ws.Cells.Item(11, 5).Value = "0.5f"   # OptionBluntMultiplier
$ws.Cells.Item(12, 5).Value = "0.6f"   # OptionFireMultiplier
$ws.Cells.Item(13, 5).Value = "1.5f"   # OptionLightningMultiplier
